# Update database: shift period columns (drop oldest FY 1396/12, add newest FY 1401/12)
# and refresh the read_price values/publish dates accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: fiscal-period headers (shift left by one period, newest goes to column H) ---
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# --- Row 9: publication dates for each period (same shift) ---
$ws.Range("D9").Value = "1399-04-18 (8)"
$ws.Range("E9").Value = "1400-04-19 (8)"
$ws.Range("F9").Value = "1401-05-19 (8)"
$ws.Range("G9").Value = "1402-02-27 (7)"
$ws.Range("H9").Value = "1402-02-27"

# --- Row 11: فروش (Sales) ---
$ws.Range("D11").Value = 26062
$ws.Range("E11").Value = 25715
$ws.Range("F11").Value = 23536
$ws.Range("G11").Value = 40149
$ws.Range("H11").Value = 57825

# --- Row 12: بهای تمام شده کالای فروش رفته (Cost of goods sold) ---
$ws.Range("D12").Value = -17596
$ws.Range("E12").Value = -20090
$ws.Range("F12").Value = -17311
$ws.Range("G12").Value = -32212
$ws.Range("H12").Value = -46116

# --- Row 13: سود (زیان) ناخالص (Gross profit) ---
$ws.Range("D13").Value = 8466
$ws.Range("E13").Value = 5624
$ws.Range("F13").Value = 6225
$ws.Range("G13").Value = 7938
$ws.Range("H13").Value = 11709

# --- Row 14: هزینه های عمومی, اداری و تشکیلاتی (G&A expenses) ---
$ws.Range("D14").Value = -855
$ws.Range("E14").Value = -1239
$ws.Range("F14").Value = -1171
$ws.Range("G14").Value = -1446
$ws.Range("H14").Value = -2791

# --- Row 15: unchanged ("-" placeholders) ---

# --- Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی (Other operating income/expense, net) ---
$ws.Range("D16").Value = 137
$ws.Range("E16").Value = -10
$ws.Range("F16").Value = 61
$ws.Range("G16").Value = 87
$ws.Range("H16").Value = 57

# --- Row 17: سود (زیان) عملیاتی (Operating profit) ---
$ws.Range("D17").Value = 7747
$ws.Range("E17").Value = 4376
$ws.Range("F17").Value = 5114
$ws.Range("G17").Value = 6579
$ws.Range("H17").Value = 8975

# --- Row 18: هزینه های مالی (Financial expenses) ---
$ws.Range("D18").Value = -328
$ws.Range("E18").Value = -221
$ws.Range("F18").Value = -92
$ws.Range("G18").Value = -149
$ws.Range("H18").Value = -1234

# --- Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی (Other non-operating income/expense, net) ---
$ws.Range("D19").Value = 215
$ws.Range("E19").Value = 208
$ws.Range("F19").Value = -2027
$ws.Range("G19").Value = 1891
$ws.Range("H19").Value = -157

# --- Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات (Pre-tax profit from continuing ops) ---
$ws.Range("D20").Value = 7634
$ws.Range("E20").Value = 4363
$ws.Range("F20").Value = 2995
$ws.Range("G20").Value = 8321
$ws.Range("H20").Value = 7584

# --- Row 21: مالیات (Tax); G21/H21 are now "-" (no tax reported yet) ---
$ws.Range("D21").Value = -8
$ws.Range("E21").Value = -31
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = "-"
$ws.Range("H21").Value = "-"

# --- Row 22: سود (زیان) خالص عملیات در حال تداوم (Profit from continuing ops) ---
$ws.Range("D22").Value = 7625
$ws.Range("E22").Value = 4332
$ws.Range("F22").Value = 2995
$ws.Range("G22").Value = 8321
$ws.Range("H22").Value = 7584

# --- Row 23: unchanged ("-" placeholders) ---

# --- Row 24: سود (زیان) خالص (Net profit) ---
$ws.Range("D24").Value = 7625
$ws.Range("E24").Value = 4332
$ws.Range("F24").Value = 2995
$ws.Range("G24").Value = 8321
$ws.Range("H24").Value = 7584

# --- Row 25: unchanged (all 0) ---

# --- Row 26: سرمایه (Capital) ---
$ws.Range("D26").Value = 1977
$ws.Range("E26").Value = 1559
$ws.Range("F26").Value = 3096
$ws.Range("G26").Value = 2653
$ws.Range("H26").Value = 1984

# --- Row 27: unchanged (all 0) ---

Write-Host "Applied database update for FY 1401/12 and shifted read_price columns."
